$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1 view: deselect its tab, move selection to E5 ---
$ws1.Range("E5").Select()

# --- Add Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Row 1
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "'Expiry (weeks)"; $arr[0,1] = "150"; $arr[0,2] = "'152.5"; $arr[0,3] = "155"; $arr[0,4] = "'157.5"; $arr[0,5] = "160"; $arr[0,6] = "'162.5"; $arr[0,7] = "165"; $arr[0,8] = "'167.5"; $arr[0,9] = "170"; $arr[0,10] = "'172.5"; $arr[0,11] = "175"; $arr[0,12] = "'177.5"; $arr[0,13] = "180"; $arr[0,14] = "'182.5"; $arr[0,15] = "185"; $arr[0,16] = "'187.5"; $arr[0,17] = "190"; $arr[0,18] = "'192.5"; $arr[0,19] = "195"; $arr[0,20] = "'197.5"; $arr[0,21] = "200"
$ws2.Range("A1:V1").Value = $arr

# Row 2
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "1"; $arr[0,1] = "'25.10"; $arr[0,2] = "'23.15"; $arr[0,3] = "'21.25"; $arr[0,4] = "'19.40"; $arr[0,5] = "'17.60"; $arr[0,6] = "'15.85"; $arr[0,7] = "'14.15"; $arr[0,8] = "'12.50"; $arr[0,9] = "'10.90"; $arr[0,10] = "'9.35"; $arr[0,11] = "'7.85"; $arr[0,12] = "'6.40"; $arr[0,13] = "'5.00"; $arr[0,14] = "'3.65"; $arr[0,15] = "'2.35"; $arr[0,16] = "'1.10"; $arr[0,17] = "'0.80"; $arr[0,18] = "'0.65"; $arr[0,19] = "'0.55"; $arr[0,20] = "'0.45"; $arr[0,21] = "'0.35"
$ws2.Range("A2:V2").Value = $arr

# Row 3
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "2"; $arr[0,1] = "'26.00"; $arr[0,2] = "'24.10"; $arr[0,3] = "'22.25"; $arr[0,4] = "'20.45"; $arr[0,5] = "'18.70"; $arr[0,6] = "'17.00"; $arr[0,7] = "'15.35"; $arr[0,8] = "'13.75"; $arr[0,9] = "'12.20"; $arr[0,10] = "'10.70"; $arr[0,11] = "'9.25"; $arr[0,12] = "'7.85"; $arr[0,13] = "'6.50"; $arr[0,14] = "'5.20"; $arr[0,15] = "'3.95"; $arr[0,16] = "'2.75"; $arr[0,17] = "'1.60"; $arr[0,18] = "'1.20"; $arr[0,19] = "'1.00"; $arr[0,20] = "'0.85"; $arr[0,21] = "'0.70"
$ws2.Range("A3:V3").Value = $arr

# Row 4
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "3"; $arr[0,1] = "'26.90"; $arr[0,2] = "'25.05"; $arr[0,3] = "'23.25"; $arr[0,4] = "'21.50"; $arr[0,5] = "'19.80"; $arr[0,6] = "'18.15"; $arr[0,7] = "'16.55"; $arr[0,8] = "'15.00"; $arr[0,9] = "'13.50"; $arr[0,10] = "'12.05"; $arr[0,11] = "'10.65"; $arr[0,12] = "'9.30"; $arr[0,13] = "'8.00"; $arr[0,14] = "'6.75"; $arr[0,15] = "'5.55"; $arr[0,16] = "'4.40"; $arr[0,17] = "'3.30"; $arr[0,18] = "'2.30"; $arr[0,19] = "'1.80"; $arr[0,20] = "'1.50"; $arr[0,21] = "'1.25"
$ws2.Range("A4:V4").Value = $arr

# Row 5
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "4"; $arr[0,1] = "'27.75"; $arr[0,2] = "'25.95"; $arr[0,3] = "'24.20"; $arr[0,4] = "'22.50"; $arr[0,5] = "'20.85"; $arr[0,6] = "'19.25"; $arr[0,7] = "'17.70"; $arr[0,8] = "'16.20"; $arr[0,9] = "'14.75"; $arr[0,10] = "'13.35"; $arr[0,11] = "'12.00"; $arr[0,12] = "'10.70"; $arr[0,13] = "'9.45"; $arr[0,14] = "'8.25"; $arr[0,15] = "'7.10"; $arr[0,16] = "'6.00"; $arr[0,17] = "'4.95"; $arr[0,18] = "'3.95"; $arr[0,19] = "'3.35"; $arr[0,20] = "'2.80"; $arr[0,21] = "'2.35"
$ws2.Range("A5:V5").Value = $arr

# Row 6
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "5"; $arr[0,1] = "'28.55"; $arr[0,2] = "'26.80"; $arr[0,3] = "'25.10"; $arr[0,4] = "'23.40"; $arr[0,5] = "'21.80"; $arr[0,6] = "'20.25"; $arr[0,7] = "'18.75"; $arr[0,8] = "'17.30"; $arr[0,9] = "'15.90"; $arr[0,10] = "'14.55"; $arr[0,11] = "'13.25"; $arr[0,12] = "'12.00"; $arr[0,13] = "'10.80"; $arr[0,14] = "'9.65"; $arr[0,15] = "'8.55"; $arr[0,16] = "'7.50"; $arr[0,17] = "'6.50"; $arr[0,18] = "'5.55"; $arr[0,19] = "'4.85"; $arr[0,20] = "'4.25"; $arr[0,21] = "'3.75"
$ws2.Range("A6:V6").Value = $arr

# Row 7
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "6"; $arr[0,1] = "'29.30"; $arr[0,2] = "'27.60"; $arr[0,3] = "'25.90"; $arr[0,4] = "'24.25"; $arr[0,5] = "'22.70"; $arr[0,6] = "'21.20"; $arr[0,7] = "'19.75"; $arr[0,8] = "'18.35"; $arr[0,9] = "'17.00"; $arr[0,10] = "'15.70"; $arr[0,11] = "'14.45"; $arr[0,12] = "'13.25"; $arr[0,13] = "'12.10"; $arr[0,14] = "'11.00"; $arr[0,15] = "'9.95"; $arr[0,16] = "'8.95"; $arr[0,17] = "'8.00"; $arr[0,18] = "'7.10"; $arr[0,19] = "'6.40"; $arr[0,20] = "'5.75"; $arr[0,21] = "'5.15"
$ws2.Range("A7:V7").Value = $arr

# Row 8
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "7"; $arr[0,1] = "'30.00"; $arr[0,2] = "'28.35"; $arr[0,3] = "'26.65"; $arr[0,4] = "'25.05"; $arr[0,5] = "'23.55"; $arr[0,6] = "'22.10"; $arr[0,7] = "'20.70"; $arr[0,8] = "'19.35"; $arr[0,9] = "'18.05"; $arr[0,10] = "'16.80"; $arr[0,11] = "'15.60"; $arr[0,12] = "'14.45"; $arr[0,13] = "'13.35"; $arr[0,14] = "'12.30"; $arr[0,15] = "'11.30"; $arr[0,16] = "'10.35"; $arr[0,17] = "'9.45"; $arr[0,18] = "'8.60"; $arr[0,19] = "'7.90"; $arr[0,20] = "'7.25"; $arr[0,21] = "'6.65"
$ws2.Range("A8:V8").Value = $arr

# Row 9
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "8"; $arr[0,1] = "'30.65"; $arr[0,2] = "'29.05"; $arr[0,3] = "'27.35"; $arr[0,4] = "'25.75"; $arr[0,5] = "'24.30"; $arr[0,6] = "'22.90"; $arr[0,7] = "'21.55"; $arr[0,8] = "'20.25"; $arr[0,9] = "'19.00"; $arr[0,10] = "'17.80"; $arr[0,11] = "'16.65"; $arr[0,12] = "'15.55"; $arr[0,13] = "'14.50"; $arr[0,14] = "'13.50"; $arr[0,15] = "'12.55"; $arr[0,16] = "'11.65"; $arr[0,17] = "'10.80"; $arr[0,18] = "'10.00"; $arr[0,19] = "'9.30"; $arr[0,20] = "'8.65"; $arr[0,21] = "'8.05"
$ws2.Range("A9:V9").Value = $arr

# Row 10
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "9"; $arr[0,1] = "'31.25"; $arr[0,2] = "'29.70"; $arr[0,3] = "'28.00"; $arr[0,4] = "'26.40"; $arr[0,5] = "'25.00"; $arr[0,6] = "'23.60"; $arr[0,7] = "'22.30"; $arr[0,8] = "'21.05"; $arr[0,9] = "'19.85"; $arr[0,10] = "'18.70"; $arr[0,11] = "'17.60"; $arr[0,12] = "'16.55"; $arr[0,13] = "'15.55"; $arr[0,14] = "'14.60"; $arr[0,15] = "'13.70"; $arr[0,16] = "'12.85"; $arr[0,17] = "'12.05"; $arr[0,18] = "'11.30"; $arr[0,19] = "'10.60"; $arr[0,20] = "'9.95"; $arr[0,21] = "'9.35"
$ws2.Range("A10:V10").Value = $arr

# Row 11
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = "10"; $arr[0,1] = "'31.80"; $arr[0,2] = "'30.30"; $arr[0,3] = "'28.60"; $arr[0,4] = "'27.00"; $arr[0,5] = "'25.65"; $arr[0,6] = "'24.30"; $arr[0,7] = "'23.00"; $arr[0,8] = "'21.80"; $arr[0,9] = "'20.65"; $arr[0,10] = "'19.55"; $arr[0,11] = "'18.50"; $arr[0,12] = "'17.50"; $arr[0,13] = "'16.55"; $arr[0,14] = "'15.65"; $arr[0,15] = "'14.80"; $arr[0,16] = "'14.00"; $arr[0,17] = "'13.25"; $arr[0,18] = "'12.55"; $arr[0,19] = "'11.90"; $arr[0,20] = "'11.30"; $arr[0,21] = "'10.75"
$ws2.Range("A11:V11").Value = $arr

# --- Style the header row (row 1) like Sheet1's header row (bold, centered, wrapped) ---
$ws1.Range("A1").Copy()
$ws2.Range("A1:V1").PasteSpecial(-4122)
$ws2.Rows.Item(1).RowHeight = 29

# --- Style the data rows (2-11) like Sheet1's data rows (wrapped, vertical center) ---
$ws1.Range("A2").Copy()
$ws2.Range("A2:V11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Sheet2 view: make it the active/selected tab, set selection to M12 ---
$ws2.Activate()
$ws2.Range("M12").Select()
